# grade update: midterm 2
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header
$ws.Range("G1").Value = "Midterm 2"

# Give column G a bit more width (matches author's manual column resize)
$ws.Columns.Item(7).ColumnWidth = 9.68

# New Midterm 2 scores (out of 60) for the students who have them
$ws.Range("G3").Formula = "=58/60"
$ws.Range("G4").Formula = "=57/60"
$ws.Range("G5").Formula = "=58/60"
$ws.Range("G9").Formula = "=60/60"
$ws.Range("G12").Formula = "=55/60"

# Match the formatting Excel applies to the Midterm 1 column when editing nearby cells
$ws.Range("F1:F12").Style = $ws.Range("A1").Style

# Move the active selection, mirroring where the author's cursor ended up
$ws.Range("G10").Select()
